# Add TODO note about auto saving and auto closing:
# append two new preference rows ("Autoclose" / "Autosave") with
# boolean FALSE values to the "Preferences" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Preferences")

$ws.Range("A7").Value = "Autoclose"
$ws.Range("B7").Value = $false

$ws.Range("A8").Value = "Autosave"
$ws.Range("B8").Value = $false

$ws.Range("A8").Select()
